# Fixed Event Update Bug
# Old Data does not stay when new data is saved.
# Status Option cannot populate when displayed, but does give validation
# warning when submitting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C11 was showing "UNKNOWN" (a status value that should no longer exist);
# restore it to the correct "Confirmed" status so old data is preserved.
$ws.Range("C11").Value = "Confirmed"

# Update the saved cursor/selection position on the sheet.
$ws.Range("C20").Select()
